$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header counts in row 1 (columns B:E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update "CON" data row (row 2, columns B:E) with new (Lichtwark-adjusted) values
$ws.Range("B2").Value = 5.0637421277820867
$ws.Range("C2").Value = 4.2841422700928815
$ws.Range("D2").Value = 6.9414908877550401
$ws.Range("E2").Value = 6.1128034389697872

# Update "STR" data row (row 3, columns B:E) with new (Lichtwark-adjusted) values
$ws.Range("B3").Value = 4.4550762181419969
$ws.Range("C3").Value = 6.8061620425162186
$ws.Range("D3").Value = 7.2646165724020548
$ws.Range("E3").Value = 5.5698631668856535

# Narrow the active selection to the edited block (B1:E3)
$ws.Range("B1:E3").Select()
